$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete row 670 entirely (shifts rows 671-856 up to become 670-855)
$ws.Rows.Item(670).Delete()
